$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: force text-number format on Price (D) cells that will receive
# plain numeric-looking strings, so Excel keeps them as text (matching the
# original inline-string cell type) instead of silently converting them to
# numbers.
$priceCells = @("D2", "D3", "D5", "D8", "D10", "D11", "D12", "D13", "D14", "D16", "D17", "D20", "D23", "D24", "D25", "D27", "D29", "D31", "D32", "D33", "D36", "D41", "D42", "D43", "D44", "D46", "D49")
foreach ($c in $priceCells) {
    $ws.Range($c).NumberFormat = "@"
}

# --- Step 2: write the updated values row by row ---

# Row 2
$ws.Range("D2").Value = "25.947.51"
$ws.Range("E2").Value = "  -0.25%  "

# Row 3
$ws.Range("D3").Value = "1.615.85"
$ws.Range("E3").Value = "  -1.14%  "

# Row 4
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").Value = "211.38"
$ws.Range("E5").Value = "  -1.25%  "

# Row 6
$ws.Range("E6").Value = "  -0.08%  "

# Row 7
$ws.Range("E7").Value = "  -3.15%  "

# Row 8
$ws.Range("D8").Value = "0.0620"
$ws.Range("E8").Value = "  -0.52%  "

# Row 9
$ws.Range("E9").Value = "  -1.53%  "

# Row 10
$ws.Range("D10").Value = "18.20"
$ws.Range("E10").Value = "  -1.48%  "

# Row 11
$ws.Range("D11").Value = "0.0792"
$ws.Range("E11").Value = "  -0.13%  "

# Row 12
$ws.Range("D12").Value = "1.840.77"
$ws.Range("E12").Value = "  -1.11%  "

# Row 15
$ws.Range("E15").Value = "  -2.00%  "

# Row 16
$ws.Range("D16").Value = "25.954.50"
$ws.Range("E16").Value = "  -0.21%  "

# Row 17
$ws.Range("D17").Value = "61.49"
$ws.Range("E17").Value = "  -0.53%  "

# Row 18
$ws.Range("E18").Value = "  -1.82%  "

# Row 19
$ws.Range("E19").Value = "  -0.10%  "

# Row 20
$ws.Range("D20").Value = "191.19"
$ws.Range("E20").Value = "  +0.36%  "

# Row 21
$ws.Range("E21").Value = "  -0.82%  "

# Row 22
$ws.Range("E22").Value = "  -1.63%  "

# Row 23
$ws.Range("D23").Value = "5.99"
$ws.Range("E23").Value = "  -2.28%  "

# Row 24
$ws.Range("D24").Value = "0.130"
$ws.Range("E24").Value = "  -2.27%  "

# Row 25
$ws.Range("D25").Value = "143.01"
$ws.Range("E25").Value = "  -0.18%  "

# Row 26
$ws.Range("E26").Value = "  -0.11%  "

# Row 27
$ws.Range("D27").Value = "1.72"
$ws.Range("E27").Value = "  -2.35%  "

# Row 28
$ws.Range("E28").Value = "  -2.19%  "

# Row 29
$ws.Range("D29").Value = "15.10"
$ws.Range("E29").Value = "  -0.73%  "

# Row 30
$ws.Range("E30").Value = "  -1.33%  "

# Row 31
$ws.Range("D31").Value = "0.0472"
$ws.Range("E31").Value = "  -2.13%  "

# Row 32
$ws.Range("D32").Value = "3.11"
$ws.Range("E32").Value = "  -1.51%  "

# Row 33
$ws.Range("D33").Value = "3.08"
$ws.Range("E33").Value = "  -2.34%  "

# Row 34
$ws.Range("E34").Value = "  -0.98%  "

# Row 35
$ws.Range("E35").Value = "  -1.14%  "

# Row 36
$ws.Range("D36").Value = "1.122.37"
$ws.Range("E36").Value = "  -0.90%  "

# Row 37
$ws.Range("E37").Value = "  -6.36%  "

# Row 38
$ws.Range("E38").Value = "  -2.21%  "

# Row 39
$ws.Range("E39").Value = "  -2.25%  "

# Row 40
$ws.Range("E40").Value = "  -1.47%  "

# Row 41
$ws.Range("D41").Value = "97.10"
$ws.Range("E41").Value = "  -1.66%  "

# Row 42
$ws.Range("D42").Value = "1.752.42"
$ws.Range("E42").Value = "  -1.09%  "

# Row 43
$ws.Range("D43").Value = "0.752"
$ws.Range("E43").Value = "  -4.11%  "

# Row 44
$ws.Range("D44").Value = "5.06"
$ws.Range("E44").Value = "  -4.48%  "

# Row 45
$ws.Range("E45").Value = "  +0.32%  "

# Row 46
$ws.Range("D46").Value = "53.73"
$ws.Range("E46").Value = "  -2.70%  "

# Row 47
$ws.Range("E47").Value = "  -0.05%  "

# Row 48
$ws.Range("E48").Value = "  -2.35%  "

# Row 49
$ws.Range("D49").Value = "0.411"
$ws.Range("E49").Value = "  -0.81%  "

# Row 50
$ws.Range("E50").Value = "  -0.12%  "

# Row 51
$ws.Range("E51").Value = "  -1.82%  "

# Row 13 <-> Row 14 swap (Polkadot and WrappedEther swapped rank position)
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.609.35"
$ws.Range("E13").Value = "  -1.44%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "4.09"
$ws.Range("E14").Value = "  -2.80%  "

# --- Step 3: restore default cell style on the Price cells we touched in Step 1
# so we do not leave a stray NumberFormat behind (matches original formatting).
foreach ($c in $priceCells) {
    $ws.Range($c).Style = "Normal"
}
